$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45177 (2023-09-08)
# to 45178 (2023-09-09) for every data row (rows 2 through 494).
$ws.Range("C2:C494").Value = 45178
